$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Extend the existing table by one row (keeps dimension/table ref/autofilter ref in sync)
$tbl = $ws.ListObjects.Item("Table2")
$tbl.ListRows.Add() | Out-Null

# Fill in the new row's values. Shared-string insertion order matters for matching
# the original author's save (dev.to link, then title, then hashnode link).
$ws.Range("B77").Value = 67
$ws.Range("F77").Value = "https://dev.to/rahulmishra05/memory-management-techniques-operating-system-m05-p02-2ei1"
$ws.Range("C77").Value = "Memory management Techniques | Operating System - M05 P02"
$ws.Range("E77").Value = "https://programmingport.hashnode.dev/memory-management-techniques-or-operating-system-m05-p02-1"
$ws.Range("D77").Value = (Get-Date -Year 2020 -Month 12 -Day 16 -Hour 0 -Minute 0 -Second 0)

# Match the author's final selection state
$ws.Range("E77").Select() | Out-Null
